$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 151, shifting existing rows 151:270 down to 152:271
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record's data
$ws.Cells.Item(151, 1).Value  = 9
$ws.Cells.Item(151, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(151, 3).Value  = "Metropolitana"
$ws.Cells.Item(151, 4).Value  = 44488
$ws.Cells.Item(151, 5).Value  = 13
$ws.Cells.Item(151, 6).Value  = "Fruta"
$ws.Cells.Item(151, 7).Value  = 100108
$ws.Cells.Item(151, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(151, 9).Value  = 100108002
$ws.Cells.Item(151, 10).Value = "Mango"
$ws.Cells.Item(151, 11).Value = "Sin especificar"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 780
$ws.Cells.Item(151, 14).Value = 6000
$ws.Cells.Item(151, 15).Value = 6500
$ws.Cells.Item(151, 16).Value = 6244
$ws.Cells.Item(151, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(151, 18).Value = "Perú"
$ws.Cells.Item(151, 19).Value = 1561
$ws.Cells.Item(151, 20).Value = 4
